# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the per-job Leve sheets, per the latest pricing pull.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1407.6923
$ws.Range("I28").Value = 1535
$ws.Range("J28").Value = 983.3333
$ws.Range("K28").Value = 1535
$ws.Range("L28").Value = 983.3333
$ws.Range("M28").Value = -1050
$ws.Range("N28").Value = -1953.3333
# Row 70
$ws.Range("H70").Value = 1889.0667
$ws.Range("I70").Value = 1286.2222
$ws.Range("J70").Value = 2793.3333
$ws.Range("K70").Value = 3858.6666
$ws.Range("L70").Value = 8379.999899999999
$ws.Range("M70").Value = -3588.6666
$ws.Range("N70").Value = -8919.999899999999
# Row 73
$ws.Range("H73").Value = 1889.0667
$ws.Range("I73").Value = 1286.2222
$ws.Range("J73").Value = 2793.3333
$ws.Range("K73").Value = 3858.6666
$ws.Range("L73").Value = 8379.999899999999
$ws.Range("M73").Value = -2922.6666
$ws.Range("N73").Value = -10251.9999
# Row 108
$ws.Range("H108").Value = 43000
$ws.Range("J108").Value = 43000
$ws.Range("L108").Value = 43000
$ws.Range("N108").Value = -50680
# Row 130
$ws.Range("H130").Value = 41853.332
$ws.Range("J130").Value = 41853.332
$ws.Range("L130").Value = 41853.332
$ws.Range("N130").Value = -51893.332
# Row 138
$ws.Range("H138").Value = 1504.9445
$ws.Range("I138").Value = 1255.5625
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 3766.6875
$ws.Range("L138").Value = 10500
$ws.Range("M138").Value = 1373.3125
$ws.Range("N138").Value = -20780
# Row 141
$ws.Range("H141").Value = 7404.05
$ws.Range("I141").Value = 9598.538
$ws.Range("J141").Value = 3328.5715
$ws.Range("K141").Value = 28795.614
$ws.Range("L141").Value = 9985.7145
$ws.Range("M141").Value = -23615.614
$ws.Range("N141").Value = -20345.7145

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3804.1567
$ws.Range("I32").Value = 3795.3064
$ws.Range("J32").Value = 3830.2856
$ws.Range("K32").Value = 3795.3064
$ws.Range("L32").Value = 3830.2856
$ws.Range("M32").Value = -3508.3064
$ws.Range("N32").Value = -4404.2856
# Row 33
$ws.Range("H33").Value = 39999.5
$ws.Range("J33").Value = 39999.5
$ws.Range("L33").Value = 39999.5
$ws.Range("N33").Value = -40657.5
# Row 36
$ws.Range("H36").Value = 9630.5
$ws.Range("I36").Value = 9260
$ws.Range("J36").Value = 10001
$ws.Range("K36").Value = 9260
$ws.Range("L36").Value = 10001
$ws.Range("M36").Value = -8914
$ws.Range("N36").Value = -10693
# Row 45
$ws.Range("H45").Value = 3541.3
$ws.Range("I45").Value = 2680
$ws.Range("J45").Value = 4402.6
$ws.Range("K45").Value = 2680
$ws.Range("L45").Value = 4402.6
$ws.Range("M45").Value = -2303
$ws.Range("N45").Value = -5156.6
# Row 137
$ws.Range("H137").Value = 45747.6
$ws.Range("J137").Value = 45747.6
$ws.Range("L137").Value = 45747.6
$ws.Range("N137").Value = -55947.6

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 95
$ws.Range("H95").Value = 32400
$ws.Range("J95").Value = 32400
$ws.Range("L95").Value = 32400
$ws.Range("N95").Value = -37892
# Row 137
$ws.Range("H137").Value = 50737.5
$ws.Range("J137").Value = 50737.5
$ws.Range("L137").Value = 50737.5
$ws.Range("N137").Value = -60937.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 216702.25
$ws.Range("I31").Value = 520599.7
$ws.Range("J31").Value = 3152.7026
$ws.Range("K31").Value = 520599.7
$ws.Range("L31").Value = 3152.7026
$ws.Range("M31").Value = -520304.7
$ws.Range("N31").Value = -3742.7026
# Row 34
$ws.Range("H34").Value = 216702.25
$ws.Range("I34").Value = 520599.7
$ws.Range("J34").Value = 3152.7026
$ws.Range("K34").Value = 520599.7
$ws.Range("L34").Value = 3152.7026
$ws.Range("M34").Value = -520397.7
$ws.Range("N34").Value = -3556.7026
# Row 141
$ws.Range("H141").Value = 26562.5
$ws.Range("J141").Value = 26562.5
$ws.Range("L141").Value = 26562.5
$ws.Range("N141").Value = -36922.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2213.85
$ws.Range("I5").Value = 1283.2307
$ws.Range("J5").Value = 3942.1428
$ws.Range("K5").Value = 3849.6921
$ws.Range("L5").Value = 11826.4284
$ws.Range("M5").Value = -3737.6921
$ws.Range("N5").Value = -12050.4284
# Row 12
$ws.Range("H12").Value = 71.76470999999999
$ws.Range("J12").Value = 106.63636
$ws.Range("L12").Value = 319.90908
$ws.Range("N12").Value = -665.90908
# Row 131
$ws.Range("H131").Value = 790.9299999999999
$ws.Range("I131").Value = 455
$ws.Range("J131").Value = 804.92706
$ws.Range("K131").Value = 1365
$ws.Range("L131").Value = 2414.78118
$ws.Range("M131").Value = 3675
$ws.Range("N131").Value = -12494.78118
# Row 135
$ws.Range("H135").Value = 2213.85
$ws.Range("I135").Value = 1283.2307
$ws.Range("J135").Value = 3942.1428
$ws.Range("K135").Value = 11549.0763
$ws.Range("L135").Value = 35479.2852
$ws.Range("M135").Value = -9014.076300000001
$ws.Range("N135").Value = -40549.2852

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
# Row 107
$ws.Range("H107").Value = 4274210.5
$ws.Range("J107").Value = 9260314
$ws.Range("L107").Value = 9260314
$ws.Range("N107").Value = -9264154
# Row 113
$ws.Range("H113").Value = 1364
$ws.Range("I113").Value = 1103.6666
$ws.Range("J113").Value = 1494.1666
$ws.Range("K113").Value = 1103.6666
$ws.Range("L113").Value = 1494.1666
$ws.Range("M113").Value = 1066.3334
$ws.Range("N113").Value = -5834.1666
# Row 137
$ws.Range("H137").Value = 39186.668
$ws.Range("J137").Value = 48780
$ws.Range("L137").Value = 48780
$ws.Range("N137").Value = -58980

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 18522332
$ws.Range("I93").Value = 27780998
$ws.Range("K93").Value = 27780998
$ws.Range("M93").Value = -27779750
# Row 132
$ws.Range("H132").Value = 3785.3333
$ws.Range("I132").Value = 2766.6667
$ws.Range("J132").Value = 5822.6665
$ws.Range("K132").Value = 8300.000100000001
$ws.Range("L132").Value = 17467.9995
$ws.Range("M132").Value = -5770.000100000001
$ws.Range("N132").Value = -22527.9995

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 889134.5600000001
$ws.Range("I126").Value = 1776
$ws.Range("J126").Value = 2663851.8
$ws.Range("K126").Value = 5328
$ws.Range("L126").Value = 7991555.399999999
$ws.Range("M126").Value = -2858
$ws.Range("N126").Value = -7996495.399999999
